$d = $word.ActiveDocument

# Real Word coalesces adjacent runs that share identical direct formatting
# once a paragraph is touched by an edit. Toggling a character property to a
# different value and back forces Word to keep (or create) an explicit run
# boundary at that spot without changing anything visible.
function Force-Split($range) {
    $range.Font.Bold = $true
    $range.Font.Bold = $false
}

# Locate the run that needs to be split: '" dan melakukan operasi pada komputer jarak jauh'
$r = $d.Content
$found = $r.Find.Execute('" dan melakukan operasi pada komputer jarak jauh', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$origStart = $r.Start
$origEnd = $r.End

$insertText = "ke sebuah sistem operasi dari jarak jauh "
$midStart = $origStart + 2
$insertRange = $d.Range($midStart, $midStart)
$insertRange.InsertBefore($insertText)

# After insertion everything from $midStart onward is shifted right by Len($insertText).
# New layout within what used to be a single run:
#   [origStart, midStart)                  -> '" '                                         (keep as its own run)
#   [midStart, midStart+40)                -> 'ke sebuah sistem operasi dari jarak jauh'    (brand-new run)
#   [midStart+40, origEnd+Len(insertText)) -> ' dan melakukan operasi pada komputer jarak jauh' (tail, keep as its own run)

$quoteRange = $d.Range($origStart, $midStart)
Force-Split $quoteRange

$newTextRange = $d.Range($midStart, $midStart + 40)
Force-Split $newTextRange

$tailStart = $midStart + 40
$tailEnd = $origEnd + $insertText.Length
$tailRange = $d.Range($tailStart, $tailEnd)
Force-Split $tailRange

# The insertion re-flows (coalesces) every run from the insertion point through the end
# of the paragraph, so the pre-existing run boundaries on either side of our edit
# ('log-in' before it, '(remote)' after it, and the later 'dat' / 'a.' split) need to be
# re-asserted explicitly, otherwise they'd silently get absorbed into neighboring runs.

$logRange = $d.Content
$null = $logRange.Find.Execute('log-in', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Force-Split $logRange

$remoteRange = $d.Content
$null = $remoteRange.Find.Execute('(remote)', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Force-Split $remoteRange

# Restore the original 'dat' / 'a.' run boundary right after '(remote)'.
$datRange = $d.Range($remoteRange.End, $remoteRange.End + 50)
Force-Split $datRange

Write-Output "Done. Inserted: [$($insertText.Trim())]"
